$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cumulative COVID death data rows (dates 2021-01-07 .. 2021-01-14)
$data = @(
    @(44203, 2788, 638, 3426),
    @(44204, 2836, 644, 3480),
    @(44205, 2918, 654, 3572),
    @(44206, 3007, 663, 3670),
    @(44207, 3102, 668, 3770),
    @(44208, 3163, 676, 3839),
    @(44209, 3260, 686, 3946),
    @(44210, 3362, 702, 4064)
)

$startRow = 84
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
}

# Update the view: scroll back to top-left and select B4 (matches author's saved view)
$ws.Range("B4").Select()
